# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E27) listed mora periods in descending
# order (2104, 2103, 2101, 2012 .. 2004). The sheet is updated so the
# periods are listed in ascending order (2004, 2005, .. 2012, 2101, 2103,
# 2104), and the "Salario Basico" column (F16:F27) is kept in sync with
# its period (period 2104 carries a Salario Basico of 41600, every other
# period carries 48000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2004", "2005", "2006", "2007", "2008", "2009", "2010", "2011", "2012", "2101", "2103", "2104")

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $period = $periods[$i]

    $ws.Range("E$row").Value = $period

    if ($period -eq "2104") {
        $ws.Range("F$row").Value = 41600
    } else {
        $ws.Range("F$row").Value = 48000
    }
}
